$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet - the whole table is being replaced / re-shaped
# (two new columns inserted: ownTeam, oppTeam; several rows added; row order changed).
$ws.Cells.Clear()

# ---- Header row ----
$headers = @("venue","date","result","ownTeam","oppTeam","batsman","totalRuns","totalBalls","total4s","total6s","sr")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# ---- Data rows ----
# Columns: venue, date, result, ownTeam, oppTeam, batsman, totalRuns, totalBalls, total4s, total6s, sr
$data = @(
    @(" Dubai (DSC)", " October 24 2020", "Kings XI won by 12 runs", "Sunrisers Hyderabad", "Kings XI Punjab", "Priyam Garg ", "3", "5", "0", "0", "60.00"),
    @(" Dubai (DSC)", " October 13 2020", "Super Kings won by 20 runs", "Sunrisers Hyderabad", "Chennai Super Kings", "Priyam Garg ", "16", "18", "1", "0", "88.88"),
    @(" Abu Dhabi", " October 18 2020", "Match tied (KKR won the one-over eliminator)", "Sunrisers Hyderabad", "Kolkata Knight Riders", "Priyam Garg ", "4", "7", "0", "0", "57.14"),
    @(" Dubai (DSC)", " September 21 2020", "RCB won by 10 runs", "Sunrisers Hyderabad", "Royal Challengers Bangalore", "Priyam Garg ", "12", "13", "1", "0", "92.30"),
    @(" Abu Dhabi", " November 06 2020", "Sunrisers won by 6 wickets (with 2 balls remaining)", "Sunrisers Hyderabad", "Royal Challengers Bangalore", "Priyam Garg ", "7", "14", "0", "0", "50.00"),
    @(" Abu Dhabi", " November 08 2020", "Capitals won by 17 runs", "Sunrisers Hyderabad", "Delhi Capitals", "Priyam Garg ", "17", "12", "0", "2", "141.66"),
    @(" Dubai (DSC)", " October 11 2020", "Royals won by 5 wickets (with 1 ball remaining)", "Sunrisers Hyderabad", "Rajasthan Royals", "Priyam Garg ", "15", "8", "1", "1", "187.50"),
    @(" Dubai (DSC)", " October 08 2020", "Sunrisers won by 69 runs", "Sunrisers Hyderabad", "Kings XI Punjab", "Priyam Garg ", "0", "1", "0", "0", "0.00"),
    @(" Dubai (DSC)", " October 02 2020", "Sunrisers won by 7 runs", "Sunrisers Hyderabad", "Chennai Super Kings", "Priyam Garg ", "51", "26", "6", "1", "196.15"),
    @(" Sharjah", " October 04 2020", "Mumbai won by 34 runs", "Sunrisers Hyderabad", "Mumbai Indians", "Priyam Garg ", "8", "7", "0", "0", "114.28")
)

# Columns that hold numeric-looking text and must stay text (matches the
# "numbers stored as text" shape of the original sheet) so the literal
# strings (e.g. "60.00") survive instead of being coerced to numbers.
$textCols = @(7, 8, 9, 10, 11)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $r + 2
    $rowData = $data[$r]
    for ($col = 1; $col -le $rowData.Length; $col++) {
        $cell = $ws.Cells.Item($rowNum, $col)
        if ($textCols -contains $col) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $rowData[$col - 1]
    }
}
